# Update "想去人数" (want-to-go count) figures in the F column on both the
# "展览" sheet and the "全部类型" sheet to reflect newly generated output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 35
$ws1.Range("F6").Value = 52
$ws1.Range("F8").Value = 3767
$ws1.Range("F10").Value = 4439
$ws1.Range("F11").Value = 490
$ws1.Range("F12").Value = 1105

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 35
$ws4.Range("F6").Value = 52
$ws4.Range("F9").Value = 3767
$ws4.Range("F11").Value = 4439
$ws4.Range("F12").Value = 490
$ws4.Range("F13").Value = 1105
